$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text interpretation for the new rows (so values like "16°C" and
# "87%" are stored literally instead of being auto-converted to numbers
# / percentages by Excel's smart input parsing).
$ws.Range("A3:C4").NumberFormat = "@"

$ws.Range("A3").Value = "2025-06-30 18:30:15"
$ws.Range("B3").Value = "16°C"
$ws.Range("C3").Value = "87%"

$ws.Range("A4").Value = "2025-06-30 18:30:47"
$ws.Range("B4").Value = "16°C"
$ws.Range("C4").Value = "87%"

# Reset the cell style back to the workbook default ("Normal") now that
# the values are safely stored as text, so no stray formatting sticks
# to the new cells.
$ws.Range("A3:C4").Style = "Normal"
